$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (rows 3-5) ---
$wsSchedule.Range("B3").Value = 46082.70833333334
$wsSchedule.Range("C3").Value = 10
$wsSchedule.Range("D3").Value = 37.8
$wsSchedule.Range("E3").Value = 345.5872875
$wsSchedule.Range("F3").Value = 9.142520833333334
$wsSchedule.Range("A4").Value = 46082.9375
$wsSchedule.Range("B4").Value = 46083.25
$wsSchedule.Range("C4").Value = 7.5
$wsSchedule.Range("D4").Value = 28.35
$wsSchedule.Range("E4").Value = 805.5930577500001
$wsSchedule.Range("F4").Value = 28.41598087301588
$wsSchedule.Range("A5").Value = 46083.41666666666
$wsSchedule.Range("B5").Value = 46083.66666666666
$wsSchedule.Range("C5").Value = 6
$wsSchedule.Range("D5").Value = 22.68
$wsSchedule.Range("E5").Value = 717.2177415
$wsSchedule.Range("F5").Value = 31.62335720899471

# --- Detailed sheet updates (rows 36-97) ---
$wsDetailed.Range("E36").Value = "OFF"
$wsDetailed.Range("B37").Value = 56.98
$wsDetailed.Range("E37").Value = "OFF"
$wsDetailed.Range("B38").Value = 57.31
$wsDetailed.Range("B39").Value = 57.36
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 57.31
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 57.06013
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 57.06
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 52.30914
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 51.50676
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 50.03655
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 52.17509
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "ON"
$wsDetailed.Range("B48").Value = 50.53328
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("E48").Value = "ON"
$wsDetailed.Range("B49").Value = 52.13164
$wsDetailed.Range("B50").Value = 51.50568
$wsDetailed.Range("B51").Value = 56.98
$wsDetailed.Range("B52").Value = 53.70908
$wsDetailed.Range("B53").Value = 56.98
$wsDetailed.Range("B54").Value = 55.22264
$wsDetailed.Range("B55").Value = 55.27768
$wsDetailed.Range("B56").Value = 55.30405
$wsDetailed.Range("B57").Value = 55.64532
$wsDetailed.Range("B58").Value = 55.62514
$wsDetailed.Range("B61").Value = 61.11969
$wsDetailed.Range("B62").Value = 65
$wsDetailed.Range("E62").Value = "OFF"
$wsDetailed.Range("B63").Value = 82.91262
$wsDetailed.Range("B64").Value = 72.01016
$wsDetailed.Range("B65").Value = 65
$wsDetailed.Range("E65").Value = "OFF"
$wsDetailed.Range("B66").Value = 65
$wsDetailed.Range("E66").Value = "OFF"
$wsDetailed.Range("B67").Value = 65
$wsDetailed.Range("E67").Value = "OFF"
$wsDetailed.Range("B68").Value = 62.19162
$wsDetailed.Range("E68").Value = "OFF"
$wsDetailed.Range("B69").Value = 65.01009999999999
$wsDetailed.Range("E69").Value = "OFF"
$wsDetailed.Range("B70").Value = 64.33967
$wsDetailed.Range("B71").Value = 63.1496
$wsDetailed.Range("B72").Value = 76.99009
$wsDetailed.Range("B73").Value = 65.01014000000001
$wsDetailed.Range("B74").Value = 57.06015
$wsDetailed.Range("B75").Value = 57.06044
$wsDetailed.Range("B76").Value = 61.40941
$wsDetailed.Range("E76").Value = "ON"
$wsDetailed.Range("B77").Value = 62.45104
$wsDetailed.Range("E77").Value = "ON"
$wsDetailed.Range("E78").Value = "ON"
$wsDetailed.Range("E79").Value = "ON"
$wsDetailed.Range("B80").Value = 57.0602
$wsDetailed.Range("E80").Value = "ON"
$wsDetailed.Range("B81").Value = 57.1172
$wsDetailed.Range("E81").Value = "ON"
$wsDetailed.Range("B82").Value = 64.89
$wsDetailed.Range("B83").Value = 66.67514
$wsDetailed.Range("B84").Value = 78
$wsDetailed.Range("B85").Value = 81.02665
$wsDetailed.Range("B86").Value = 84.79000000000001
$wsDetailed.Range("B87").Value = 84.79000000000001
$wsDetailed.Range("B88").Value = 84.79000000000001
$wsDetailed.Range("B89").Value = 82.35193
$wsDetailed.Range("B90").Value = 83.64431
$wsDetailed.Range("B91").Value = 79.98466000000001
$wsDetailed.Range("B92").Value = 79.28679
$wsDetailed.Range("B93").Value = 78
$wsDetailed.Range("B95").Value = 58.97834
$wsDetailed.Range("B96").Value = 59.25448
$wsDetailed.Range("B97").Value = 57.31

Write-Host "Applied all updates (run 206)"
